# Applies updated "Inscritos" / "Pagos" / "Inscrições homologadas" counts
# to the Inscricoes worksheet, matching the commit's data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Inscricoes")

$updates = @(
    @{ Row = 4;  Col = "E"; Value = 51 }

    @{ Row = 5;  Col = "E"; Value = 159 }

    @{ Row = 10; Col = "E"; Value = 682 }
    @{ Row = 10; Col = "F"; Value = 379 }
    @{ Row = 10; Col = "H"; Value = 474 }

    @{ Row = 11; Col = "E"; Value = 451 }
    @{ Row = 11; Col = "F"; Value = 254 }
    @{ Row = 11; Col = "H"; Value = 319 }

    @{ Row = 12; Col = "E"; Value = 683 }

    @{ Row = 15; Col = "E"; Value = 198 }

    @{ Row = 16; Col = "E"; Value = 233 }

    @{ Row = 17; Col = "E"; Value = 125 }

    @{ Row = 23; Col = "E"; Value = 224 }
    @{ Row = 23; Col = "F"; Value = 116 }
    @{ Row = 23; Col = "H"; Value = 168 }

    @{ Row = 24; Col = "E"; Value = 263 }
    @{ Row = 24; Col = "F"; Value = 155 }
    @{ Row = 24; Col = "H"; Value = 185 }

    @{ Row = 25; Col = "E"; Value = 326 }

    @{ Row = 26; Col = "E"; Value = 190 }
    @{ Row = 26; Col = "F"; Value = 120 }
    @{ Row = 26; Col = "H"; Value = 145 }

    @{ Row = 27; Col = "F"; Value = 207 }
    @{ Row = 27; Col = "H"; Value = 289 }

    @{ Row = 28; Col = "E"; Value = 224 }

    @{ Row = 29; Col = "E"; Value = 193 }
    @{ Row = 29; Col = "F"; Value = 117 }
    @{ Row = 29; Col = "H"; Value = 158 }

    @{ Row = 30; Col = "E"; Value = 250 }
    @{ Row = 30; Col = "F"; Value = 156 }
    @{ Row = 30; Col = "H"; Value = 208 }

    @{ Row = 31; Col = "E"; Value = 81 }
    @{ Row = 31; Col = "F"; Value = 36 }
    @{ Row = 31; Col = "H"; Value = 64 }

    @{ Row = 32; Col = "E"; Value = 212 }

    @{ Row = 33; Col = "E"; Value = 326 }

    @{ Row = 34; Col = "E"; Value = 247 }
    @{ Row = 34; Col = "F"; Value = 175 }
    @{ Row = 34; Col = "H"; Value = 213 }

    @{ Row = 35; Col = "E"; Value = 181 }
    @{ Row = 35; Col = "F"; Value = 126 }
    @{ Row = 35; Col = "H"; Value = 153 }

    @{ Row = 37; Col = "E"; Value = 194 }
    @{ Row = 37; Col = "F"; Value = 109 }
    @{ Row = 37; Col = "H"; Value = 146 }

    @{ Row = 38; Col = "E"; Value = 106 }
    @{ Row = 38; Col = "F"; Value = 64 }
    @{ Row = 38; Col = "H"; Value = 81 }

    @{ Row = 39; Col = "E"; Value = 197 }

    @{ Row = 40; Col = "E"; Value = 304 }

    @{ Row = 41; Col = "E"; Value = 433 }
    @{ Row = 41; Col = "F"; Value = 220 }
    @{ Row = 41; Col = "H"; Value = 312 }

    @{ Row = 42; Col = "E"; Value = 455 }
    @{ Row = 42; Col = "F"; Value = 262 }
    @{ Row = 42; Col = "H"; Value = 323 }

    @{ Row = 43; Col = "E"; Value = 140 }
    @{ Row = 43; Col = "F"; Value = 77 }
    @{ Row = 43; Col = "H"; Value = 104 }

    @{ Row = 44; Col = "E"; Value = 361 }

    @{ Row = 45; Col = "E"; Value = 178 }

    @{ Row = 46; Col = "E"; Value = 387 }

    @{ Row = 47; Col = "E"; Value = 534 }
    @{ Row = 47; Col = "F"; Value = 302 }
    @{ Row = 47; Col = "H"; Value = 394 }

    @{ Row = 48; Col = "E"; Value = 265 }
    @{ Row = 48; Col = "F"; Value = 128 }
    @{ Row = 48; Col = "H"; Value = 172 }

    @{ Row = 49; Col = "E"; Value = 337 }

    @{ Row = 50; Col = "E"; Value = 275 }
    @{ Row = 50; Col = "F"; Value = 151 }
    @{ Row = 50; Col = "H"; Value = 224 }

    @{ Row = 51; Col = "F"; Value = 131 }
    @{ Row = 51; Col = "H"; Value = 205 }
)

foreach ($u in $updates) {
    $ws.Range("$($u.Col)$($u.Row)").Value = $u.Value
}
